# Update the quarterly income-statement "Overview" sheet:
#  - drop the oldest quarter column (D) and append the newest quarter (new M)
#  - this shifts every period header, publish-date header and data row one
#    column to the left (D<-E, E<-F, ... L<-M) and fills the new M with the
#    freshly reported quarter's values (also corrects a couple of historical
#    figures per the updated read_price algorithm)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: period-end captions (دوره مالی) ---
$row8 = @(
    "12 ماهه منتهی به 1399/09",
    "3 ماهه منتهی به 1399/12",
    "6 ماهه منتهی به 1400/03",
    "9 ماهه منتهی به 1400/06",
    "12 ماهه منتهی به 1400/09",
    "3 ماهه منتهی به 1400/12",
    "6 ماهه منتهی به 1401/03",
    "9 ماهه منتهی به 1401/06",
    "12 ماهه منتهی به 1401/09",
    "3 ماهه منتهی به 1401/12"
)

# --- Row 9: publish-date captions (تاریخ انتشار) ---
# (the bare "yyyy-mm-dd" captions need a leading apostrophe so Excel keeps
#  them as text instead of parsing them as a date serial)
$row9 = @(
    "1400-12-15 (8)",
    "1401-01-30 (2)",
    "1401-06-13 (4)",
    "1401-07-30 (2)",
    "1401-12-06 (8)",
    "1402-01-30 (2)",
    "1401-06-13 (2)",
    "'1401-07-30",
    "1402-01-30 (3)",
    "'1402-01-30"
)

for ($i = 0; $i -lt 10; $i++) {
    $col = 4 + $i
    $ws.Cells.Item(8, $col).Value = $row8[$i]
    $ws.Cells.Item(9, $col).Value = $row9[$i]
}

# --- Data rows: each row holds 10 quarters in columns D..M ---
$data = @{
    11 = @(4840393, 1727035, 4064257, 6325618, 8904537, 2632695, 5646530, 9897037, 14085615, 3139808)
    12 = @(-2399151, -965829, -2034016, -3083859, -4234148, -1375192, -2571303, -4762552, -6412052, -1637053)
    13 = @(2441242, 761206, 2030241, 3241759, 4670389, 1257503, 3075227, 5134485, 7673563, 1502755)
    14 = @(-354464, -65436, -234476, -314977, -460810, -100506, -277944, -429132, -738874, -135259)
    16 = @(-48883, 25040, -8616, -22313, -32379, 0, 17827, 17827, 304971, 0)
    17 = @(2037895, 720810, 1787149, 2904469, 4177200, 1156997, 2815110, 4723180, 7239660, 1367496)
    18 = @(-67097, -14624, -19444, -24278, -26594, 0, -625, -625, -5279, 0)
    19 = @(245788, -8440, 6203, 164249, 272876, 0, 318717, 383151, 580675, 307792)
    20 = @(2216586, 697746, 1773908, 3044440, 4423482, 1156997, 3133202, 5105706, 7815056, 1675288)
    21 = @(-285700, -134579, -288096, -528954, -556809, -220822, -391021, -815356, -674720, -177940)
    22 = @(1930886, 563167, 1485812, 2515486, 3866673, 936175, 2742181, 4290350, 7140336, 1497348)
    24 = @(1930886, 563167, 1485812, 2515486, 3866673, 936175, 2742181, 4290350, 7140336, 1497348)
    25 = @(1175, 343, 603, 1530, 1568, 380, 1112, 1740, 2896, 607)
    26 = @(1644000, 1644000, 2466000, 1644000, 2466000, 2466000, 2466000, 2466000, 2466000, 2466000)
    27 = @(783, 228, 603, 1020, 1568, 380, 1112, 1740, 2896, 607)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt 10; $i++) {
        $col = 4 + $i
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
